$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-262:
# every cell currently holding serial date 45177 is bumped to 45178 (one day later).
$range = $ws.Range("C2:C262")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
